$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory (column H) text labels for specific rows
$ws.Range("H4").Value  = "scatter plot(s)"
$ws.Range("H17").Value = "photo(s)"
$ws.Range("H20").Value = "data display"
$ws.Range("H21").Value = "photo(s)"
$ws.Range("H23").Value = "photo(s)"
$ws.Range("H27").Value = "photo(s)"
$ws.Range("H28").Value = "photo(s)"
$ws.Range("H29").Value = "data display"
$ws.Range("H30").Value = "mixed statistical plot (more than 1 statistical plot and type)"
$ws.Range("H33").Value = "drawing(s)"
$ws.Range("H40").Value = "photo(s)"
$ws.Range("H44").Value = "data display"
$ws.Range("H46").Value = "photo(s)"
$ws.Range("H47").Value = "photo(s)"
$ws.Range("H48").Value = "photo(s)"
$ws.Range("H54").Value = "drawing(s)"
$ws.Range("H55").Value = "drawing(s)"
$ws.Range("H57").Value = "photo(s)"
$ws.Range("H60").Value = "photo(s)"
$ws.Range("H71").Value = "data collection, data analysis, data gathering diagram"

# Remove the now-unneeded "is_viewed" column (column I)
$ws.Columns.Item(9).Delete()
